$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/practitioner-hierarchy-level-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")
# Row 2 is the root "Extension" element; its Constraint(s) column (AI) text
# (the ele-1/ext-1 constraint) is now expressed further down the table on the
# rows that actually define Extension.id / Extension.extension, so clear it here.
$elements.Range("AI2").Value = ""
# Row 5 (Extension.url) Fixed Value mirrors the structure definition URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/practitioner-hierarchy-level-code"
